$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column values that look numeric (e.g. "559.60") stay as text,
# matching the original inlineStr cell type in the workbook.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '61.936.62'
$ws.Range("E2").Value = '  +2.83%  '

$ws.Range("D3").Value = '2.408.39'
$ws.Range("E3").Value = '  +3.83%  '

$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '559.60'
$ws.Range("E5").Value = '  +2.66%  '

$ws.Range("D6").Value = '138.32'
$ws.Range("E6").Value = '  +5.58%  '

$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("D8").Value = '0.585'
$ws.Range("E8").Value = '  +0.86%  '

$ws.Range("D9").Value = '2.406.35'
$ws.Range("E9").Value = '  +3.87%  '

$ws.Range("E10").Value = '  +3.34%  '

$ws.Range("E11").Value = '  +3.87%  '

$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("E13").Value = '  +4.02%  '

$ws.Range("D14").Value = '25.68'
$ws.Range("E14").Value = '  +8.56%  '

$ws.Range("D15").Value = '2.835.52'
$ws.Range("E15").Value = '  +3.76%  '

$ws.Range("D16").Value = '61.895.40'
$ws.Range("E16").Value = '  +2.82%  '

$ws.Range("D17").Value = '0.0000139'
$ws.Range("E17").Value = '  +5.10%  '

$ws.Range("D18").Value = '2.420.17'
$ws.Range("E18").Value = '  +4.76%  '

$ws.Range("E19").Value = '  +4.44%  '

$ws.Range("D20").Value = '343.69'
$ws.Range("E20").Value = '  +9.52%  '

$ws.Range("E21").Value = '  +2.03%  '

$ws.Range("D22").Value = '6.84'
$ws.Range("E22").Value = '  +3.10%  '

$ws.Range("E23").Value = '  +0.35%  '

$ws.Range("D24").Value = '64.98'
$ws.Range("E24").Value = '  +2.08%  '

$ws.Range("E25").Value = '  +0.33%  '

$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '8.33'
$ws.Range("E27").Value = '  +6.17%  '

$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '1.52'
$ws.Range("E28").Value = '  +13.03%  '

$ws.Range("E29").Value = '  +15.67%  '

$ws.Range("E30").Value = '  +4.11%  '

$ws.Range("D31").Value = '0.0₃0779'
$ws.Range("E31").Value = '  +6.88%  '

$ws.Range("D32").Value = '6.36'
$ws.Range("E32").Value = '  +7.31%  '

$ws.Range("D33").Value = '170.83'
$ws.Range("E33").Value = '  -1.46%  '

$ws.Range("D34").Value = '0.397'
$ws.Range("E34").Value = '  +4.58%  '

$ws.Range("B35").Value = 'Bittensor'
$ws.Range("C35").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D35").Value = '378.98'
$ws.Range("E35").Value = '  +17.37%  '

$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '1.40'
$ws.Range("E36").Value = '  +2.44%  '

$ws.Range("D37").Value = '18.52'
$ws.Range("E37").Value = '  +3.95%  '

$ws.Range("D38").Value = '4.51'
$ws.Range("E38").Value = '  +11.45%  '

$ws.Range("E40").Value = '  -0.01%  '

$ws.Range("E41").Value = '  +9.10%  '

$ws.Range("D42").Value = '39.11'
$ws.Range("E42").Value = '  +3.15%  '

$ws.Range("D43").Value = '144.67'
$ws.Range("E43").Value = '  +4.64%  '

$ws.Range("D44").Value = '3.65'
$ws.Range("E44").Value = '  +4.81%  '

$ws.Range("D45").Value = '20.61'
$ws.Range("E45").Value = '  +7.59%  '

$ws.Range("D46").Value = '0.0528'
$ws.Range("E46").Value = '  +6.60%  '

$ws.Range("D47").Value = '0.0959'
$ws.Range("E47").Value = '  +2.05%  '

$ws.Range("D48").Value = '0.585'
$ws.Range("E48").Value = '  +4.77%  '

$ws.Range("D49").Value = '17.92'
$ws.Range("E49").Value = '  +6.16%  '

$ws.Range("D50").Value = '0.0220'
$ws.Range("E50").Value = '  +3.32%  '

$ws.Range("D51").Value = '0.0₆0216'
$ws.Range("E51").Value = '  +1.19%  '
